# Adds a new forecast column (AC, "2020-04-30") and a new forecast row (41,
# "as of" 2020-05-14) to both the "cases" and "deaths" sheets, plus fills in
# the newly-observed value for 2020-04-30 in column B.
#
# xlPasteValues = -4163 (used so the scratch cell's *value* is copied without
# carrying along any number-format/style, and so that a date-look-alike
# string like "2020-04-30" lands as literal text instead of being
# re-interpreted as a date serial number).

$wb = $excel.ActiveWorkbook

$sheetNames = @("cases", "deaths")

# New "as of" row values, keyed by sheet name, column letter -> value.
# Row 27 (as-of 2020-04-30) previously had an empty "Observed" cell (col B);
# it is now populated for both sheets.
$observedRow27 = @{ "cases" = 85380; "deaths" = 5901 }

# New column AC (the "2020-04-30" forecast column), rows 28..41, per sheet.
$newColCases = @{
    28 = 93111
    29 = 101323
    30 = 110241
    31 = 119803
    32 = 130264
    33 = 141690
    34 = 152338
    35 = 162414
    36 = 170960
    37 = 179848
    38 = 188530
    39 = 197634
    40 = 204842
    41 = 211978
}

$newColDeaths = @{
    28 = 6498
    29 = 7132
    30 = 7690
    31 = 8430
    32 = 9144
    33 = 9851
    34 = 10531
    35 = 11161
    36 = 11700
    37 = 12243
    38 = 12761
    39 = 13227
    40 = 13665
    41 = 14076
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # --- New column header AC1: reuse the existing "2020-04-30" text (it is
    # already a shared string used by A27) without triggering Excel's
    # automatic date-parsing of typed-in date-like strings. We stage the
    # literal text in an out-of-the-way scratch cell via a formula (so it is
    # never "typed" as a value), copy it, and paste-special *values only*
    # into AC1; this keeps AC1 a plain, unstyled text cell.
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = "=""2020-04-30"""
    $scratch.Copy() | Out-Null
    $ws.Range("AC1").PasteSpecial(-4163) | Out-Null
    $scratch.Clear() | Out-Null

    # --- Newly observed value for the "as of 2020-04-30" row (row 27),
    # column B ("Observed").
    $ws.Range("B27").Value = $observedRow27[$name]

    # --- New diagonal forecast column AC, rows 28-41.
    $col = $null
    if ($name -eq "cases") { $col = $newColCases } else { $col = $newColDeaths }
    foreach ($r in 28..41) {
        $ws.Cells.Item($r, 29).Value = $col[$r]
    }

    # --- New row 41: "as of 2020-05-14" (brand-new date, so a new shared
    # string entry is required). Same scratch-cell trick as above so the
    # text lands as a plain, unstyled text cell instead of a date serial.
    $scratch2 = $ws.Range("ZZ2")
    $scratch2.Formula = "=""2020-05-14"""
    $scratch2.Copy() | Out-Null
    $ws.Range("A41").PasteSpecial(-4163) | Out-Null
    $scratch2.Clear() | Out-Null
}
